$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This sheet logs one weekly price observation per row for "Betarraga"
# (rows 295..361). A new weekly record was inserted at row 295, pushing
# every existing record at row 295..361 down by one row (new row 362
# absorbs what used to be row 361). The columns that actually vary
# between observations are: D (Fecha), I (Calidad), J (Volumen),
# K/L/M (Precio minimo/maximo/promedio), O (Origen), P (Precio $/Kg).
# Everything else (A,B,C,E,F,G,H,N,Q,R) is constant across every row in
# this block, so the newly created row 362 can just copy it straight
# from row 361.

$varCols = @(4, 9, 10, 11, 12, 13, 15, 16)   # D, I, J, K, L, M, O, P
$fixedCols = @(1, 2, 3, 5, 6, 7, 8, 14, 17, 18)  # A, B, C, E, F, G, H, N, Q, R

$firstRow = 295
$lastRowOld = 361
$lastRowNew = 362

# 1) Populate the brand-new row (362) with the constant columns copied
#    from row 361 (same market/category/etc. for this whole block).
foreach ($c in $fixedCols) {
    $v = $ws.Cells.Item($lastRowOld, $c).Value()
    $ws.Cells.Item($lastRowNew, $c).Value = $v
}

# Match the date cell's number format so the new row's style matches
# the rest of the "Fecha" column.
$dateFmt = $ws.Cells.Item($lastRowOld, 4).NumberFormat()
$ws.Cells.Item($lastRowNew, 4).NumberFormat = $dateFmt

# 2) Shift the variable columns down by one row, starting from the
#    bottom so we never clobber a value before it has been copied.
for ($n = $lastRowNew; $n -ge ($firstRow + 1); $n--) {
    foreach ($c in $varCols) {
        $v = $ws.Cells.Item($n - 1, $c).Value()
        $ws.Cells.Item($n, $c).Value = $v
    }
}

# 3) Write the new weekly observation into row 295 (Volumen/Calidad keep
#    the values already there; only the date and the three price columns
#    plus the per-kg price change).
$ws.Cells.Item($firstRow, 4).Value = 44782   # Fecha
$ws.Cells.Item($firstRow, 11).Value = 750    # Precio minimo
$ws.Cells.Item($firstRow, 12).Value = 750    # Precio maximo
$ws.Cells.Item($firstRow, 13).Value = 750    # Precio promedio ponderado
$ws.Cells.Item($firstRow, 16).Value = 150    # Precio $/Kg
